$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Copy the cell formatting from existing rows so the new rows inherit the
# same styles used elsewhere in the sheet (font Arial 10 / hyperlink style /
# wrap-text style), instead of the workbook's default font.
# Row 7 is a plain data row (A-D plain text, F plain text, G hyperlink,
# H phone, I/J plain) - good template for rows 23 and 25.
# Row 11 has a wrapped/long description in F (style s="3") - good template
# for F24.
# ---------------------------------------------------------------------------

$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A23:D23").PasteSpecial(-4122) | Out-Null
$ws.Range("F7:G7").Copy() | Out-Null
$ws.Range("F23:G23").PasteSpecial(-4122) | Out-Null

$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A25:D25").PasteSpecial(-4122) | Out-Null
$ws.Range("F7:G7").Copy() | Out-Null
$ws.Range("F25:G25").PasteSpecial(-4122) | Out-Null

$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A24:D24").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Copy() | Out-Null
$ws.Range("F24").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G24").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# H column (phone numbers) should be plain numbers using the same font as
# the rest of the row (style s="1"), not the "@"-text style used by the
# other phone cells in this sheet - copy format from I23 (already s="1").
$ws.Range("I23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null
$ws.Range("I23").Copy() | Out-Null
$ws.Range("H24").PasteSpecial(-4122) | Out-Null
$ws.Range("I23").Copy() | Out-Null
$ws.Range("H25").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 23 - יאיר וסקר / לי-בר פירות
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "יאיר וסקר"
$ws.Range("B23").Value = "Hen0411@gmail.com"
$ws.Range("C23").Value = "לי-בר פירות"
$ws.Range("D23").Value = "אוכל"
$ws.Range("F23").Value = "מגשי פירות מעוצבים לכל סוג של אירוע. החל מ150 שח"
$ws.Range("G23").Value = "https://www.facebook.com/lybr.pyrwt?mibextid=LQQJ4d"
$ws.Range("H23").Value = 506323751

# ---------------------------------------------------------------------------
# Row 24 - אבי קלבו / ידידים - סיוע בדרכים
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "אבי קלבו"
$ws.Range("B24").Value = "kalvo007@gmail.com"
$ws.Range("C24").Value = "ידידים - סיוע בדרכים"
$ws.Range("D24").Value = "סיוע בדרכים"
$ws.Range("F24").Value = "הארגון מספק עזרה ראשונה לא רפואית ללא כל עלות, בדרכים ובבתים בשלל תחומים, בהם: הנעת הרכב, סיוע בהחלפת גלגל, פתיחת רכב שננעל ועוד.`nהארגון, הפועל בהתנדבות מלאה, מספק עזרה לזולת 24 שעות ביממה, בכל ימות השבוע, למעט שבתות וחגים."
$ws.Range("G24").Value = "https://yedidim-il.org/"
$ws.Range("H24").Value = 1230
$ws.Rows.Item(24).RowHeight = 294

# ---------------------------------------------------------------------------
# Row 25 - יניר שוקרון / סטודיו בייגלה
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "יניר שוקרון"
$ws.Range("B25").Value = "yanirshukrun@gmail.com"
$ws.Range("C25").Value = "סטודיו בייגלה"
$ws.Range("D25").Value = "התחום שלי לא מופיע"
$ws.Range("F25").Value = "הסטודיו מתמקד ומתקמצע בבניית אתרים על פלטפורמת וורדפרס, תחזוקה לאתרים ושיווק דיגיטלי לעסקים קטנים ובינוניים."
$ws.Range("G25").Value = "beigale.co.il"
$ws.Range("H25").Value = 542009876

# ---------------------------------------------------------------------------
# Hyperlinks for the new G cells (G23/G24 display text == address, so no
# explicit display text needed; G25 shows "beigale.co.il" while the
# resolved address is "http://beigale.co.il/" - add the hyperlink with the
# address as the display text first (so the workbook records the resolved
# address), then restore the real cell text afterwards).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G23"), "https://www.facebook.com/lybr.pyrwt?mibextid=LQQJ4d") | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("G23").Value = "https://www.facebook.com/lybr.pyrwt?mibextid=LQQJ4d"

$ws.Hyperlinks.Add($ws.Range("G24"), "https://yedidim-il.org/") | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G24").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("G24").Value = "https://yedidim-il.org/"

$ws.Hyperlinks.Add($ws.Range("G25"), "http://beigale.co.il/", "", "", "http://beigale.co.il/") | Out-Null
$ws.Range("G7").Copy() | Out-Null
$ws.Range("G25").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("G25").Value = "beigale.co.il"

# ---------------------------------------------------------------------------
# View state: scroll down and select J27 (best-effort; the runtime does not
# surface a window/topLeftCell object, so only the selection can be set).
# ---------------------------------------------------------------------------
$ws.Range("J27").Select()
